# Update "想去人数" (column F) figures across sheets to reflect refreshed
# source data (gh-pages output regenerated at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 935
$wsExhibition.Range("F3").Value = 1031
$wsExhibition.Range("F4").Value = 821
$wsExhibition.Range("F5").Value = 890
$wsExhibition.Range("F6").Value = 471
$wsExhibition.Range("F7").Value = 727
$wsExhibition.Range("F8").Value = 172
$wsExhibition.Range("F9").Value = 1325
$wsExhibition.Range("F11").Value = 429
$wsExhibition.Range("F12").Value = 572
$wsExhibition.Range("F14").Value = 62
$wsExhibition.Range("F15").Value = 62
$wsExhibition.Range("F16").Value = 1275
$wsExhibition.Range("F17").Value = 151
$wsExhibition.Range("F19").Value = 435
$wsExhibition.Range("F22").Value = 609
$wsExhibition.Range("F24").Value = 668
$wsExhibition.Range("F26").Value = 1124
$wsExhibition.Range("F27").Value = 5

# Sheet "演出" (Performance)
$wsPerformance = $wb.Worksheets.Item("演出")
$wsPerformance.Range("F9").Value = 33

# Sheet "全部类型" (All Types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 935
$wsAll.Range("F5").Value = 1031
$wsAll.Range("F6").Value = 821
$wsAll.Range("F7").Value = 890
$wsAll.Range("F8").Value = 471
$wsAll.Range("F9").Value = 727
$wsAll.Range("F10").Value = 172
$wsAll.Range("F11").Value = 1325
$wsAll.Range("F15").Value = 429
$wsAll.Range("F16").Value = 572
$wsAll.Range("F19").Value = 62
$wsAll.Range("F20").Value = 62
$wsAll.Range("F21").Value = 1275
$wsAll.Range("F23").Value = 151
$wsAll.Range("F25").Value = 435
$wsAll.Range("F30").Value = 609
$wsAll.Range("F31").Value = 33
$wsAll.Range("F36").Value = 668
$wsAll.Range("F38").Value = 1124
$wsAll.Range("F39").Value = 5
